$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update ExisUnits (column F) for newly added existing VRES/BESS units
$ws.Range("F8").Value = 7
$ws.Range("F12").Value = 16
$ws.Range("F14").Value = 4
$ws.Range("F15").Value = 28
$ws.Range("F16").Value = 47
$ws.Range("F17").Value = 8
$ws.Range("F18").Value = 2

# Update MaxInvest (column I) to the new uniform value
$ws.Range("I8").Value = 8
$ws.Range("I9").Value = 8
$ws.Range("I10").Value = 8
$ws.Range("I11").Value = 8
$ws.Range("I12").Value = 8
$ws.Range("I13").Value = 8
$ws.Range("I14").Value = 8
$ws.Range("I15").Value = 8
$ws.Range("I16").Value = 8
$ws.Range("I17").Value = 8
$ws.Range("I18").Value = 8

# Update the active cell selection to match the final state recorded in the sheet view
$ws.Range("I24").Select()
